# Update "gh-pages" generated output workbook with the latest scraped data.
# - Bumps a handful of "want-to-go" counts (column F) that ticked up/down
#   since the previous scrape, on both the "展览" (exhibitions) sheet and
#   the "全部类型" (all types) roll-up sheet.
# - Inserts a newly-scraped event ("张家港·META萌圆饿了") as a new row,
#   immediately above the existing last row ("苏州·星部落动漫嘉年华"),
#   on both of those same two sheets.

$wb = $excel.ActiveWorkbook

function Update-WantCounts {
    param($ws, $rowMap)
    foreach ($row in $rowMap.Keys) {
        $ws.Range("F" + $row).Value = $rowMap[$row]
    }
}

function Insert-NewEvent {
    param($ws, [int]$insertAt)

    # Push everything from $insertAt downward by one row, opening up a
    # fresh, completely blank row at $insertAt.
    $ws.Rows.Item($insertAt).Insert()

    $newRow = $insertAt
    $shiftedRow = $insertAt + 1

    # --- A: sequence number. Keep the same numeric/bordered/bold style the
    # rest of column A uses (the blank row Insert() leaves behind loses the
    # border, so explicitly restore it before writing the value). ---
    $aCell = $ws.Range("A" + $newRow)
    $aCell.Font.Bold = $true
    $aCell.Borders.LineStyle = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Value = ($insertAt - 1)

    # --- B: start date, stored as plain text (e.g. "2024-11-16"), not as
    # an auto-converted date serial. Force text formatting for the write,
    # then drop back to the sheet's normal style so no stray number format
    # is left behind on the cell. ---
    $bCell = $ws.Range("B" + $newRow)
    $bCell.NumberFormat = "@"
    $bCell.Value = "2024-11-16"
    $bCell.Style = "Normal"

    $ws.Range("C" + $newRow).Value = "张家港·META萌圆饿了"
    $ws.Range("D" + $newRow).Value = "杨舍镇人民中路42号 张家港国贸酒店"
    $ws.Range("E" + $newRow).Value = "2024.11.16 10:00-11.16 17:00"
    $ws.Range("F" + $newRow).Value = 15
    $ws.Range("G" + $newRow).Value = 40
    $ws.Range("H" + $newRow).Value = "https://show.bilibili.com/platform/detail.html?id=90745"
    $ws.Range("I" + $newRow).Value = "//i2.hdslb.com/bfs/openplatform/202408/jB7b2kZ11723621730632.jpeg"

    # The row that got pushed down keeps all of its own data (Insert()
    # shifts it intact) except its running sequence number in column A,
    # which advances by one to account for the newly inserted row above it.
    $shiftedACell = $ws.Range("A" + $shiftedRow)
    $shiftedACell.Value = $insertAt
}

# ---- Sheet "展览" (exhibitions): 23 rows of data -> 24 ----
$ws1 = $wb.Worksheets.Item("展览")

Update-WantCounts $ws1 @{
    3  = 1210
    4  = 16876
    5  = 34
    6  = 1646
    9  = 390
    10 = 224
    12 = 11700
    14 = 1368
    15 = 4646
    16 = 460
    17 = 2
    21 = 340
}

Insert-NewEvent $ws1 23

# ---- Sheet "全部类型" (all types roll-up): 26 rows of data -> 27 ----
$ws4 = $wb.Worksheets.Item("全部类型")

Update-WantCounts $ws4 @{
    4  = 1210
    5  = 16876
    6  = 34
    7  = 1646
    10 = 390
    11 = 224
    15 = 11700
    17 = 1368
    18 = 4646
    19 = 460
    20 = 2
    24 = 340
}

Insert-NewEvent $ws4 26
